$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update training_history_run_6 values (rows 2-51) per new training run
# (128 dense layers... 20 lstm... 50 epochs... 0.5 dropout)

$ws.Cells.Item(2, 1).Value = 0.06021292507648468
$ws.Cells.Item(2, 2).Value = 0.984203040599823
$ws.Cells.Item(2, 3).Value = 0.01489775814116001
$ws.Cells.Item(2, 4).Value = 0.9976165890693665
$ws.Cells.Item(3, 1).Value = 0.009292903356254101
$ws.Cells.Item(3, 2).Value = 0.9983160495758057
$ws.Cells.Item(3, 3).Value = 0.01305338460952044
$ws.Cells.Item(3, 4).Value = 0.9976165890693665
$ws.Cells.Item(4, 1).Value = 0.004931285511702299
$ws.Cells.Item(4, 2).Value = 0.9985566139221191
$ws.Cells.Item(4, 3).Value = 0.009170049801468849
$ws.Cells.Item(4, 4).Value = 0.9983614087104797
$ws.Cells.Item(5, 1).Value = 0.002203599316999316
$ws.Cells.Item(5, 2).Value = 0.9994387030601501
$ws.Cells.Item(5, 3).Value = 0.006220538634806871
$ws.Cells.Item(5, 4).Value = 0.9986593127250671
$ws.Cells.Item(6, 1).Value = 0.001202365849167109
$ws.Cells.Item(6, 2).Value = 0.9997193217277527
$ws.Cells.Item(6, 3).Value = 0.005590631160885096
$ws.Cells.Item(6, 4).Value = 0.9988082647323608
$ws.Cells.Item(7, 1).Value = 0.001036283327266574
$ws.Cells.Item(7, 2).Value = 0.9996591806411743
$ws.Cells.Item(7, 3).Value = 0.007719030603766441
$ws.Cells.Item(7, 4).Value = 0.9992551803588867
$ws.Cells.Item(8, 1).Value = 0.00116246216930449
$ws.Cells.Item(8, 2).Value = 0.999739408493042
$ws.Cells.Item(8, 3).Value = 0.005785822402685881
$ws.Cells.Item(8, 4).Value = 0.999106228351593
$ws.Cells.Item(9, 1).Value = 0.0006684120162390172
$ws.Cells.Item(9, 2).Value = 0.9998797178268433
$ws.Cells.Item(9, 3).Value = 0.002921943552792072
$ws.Cells.Item(9, 4).Value = 0.9992551803588867
$ws.Cells.Item(10, 1).Value = 0.0008140974678099155
$ws.Cells.Item(10, 2).Value = 0.9997594356536865
$ws.Cells.Item(10, 3).Value = 0.00329586723819375
$ws.Cells.Item(10, 4).Value = 0.9994041323661804
$ws.Cells.Item(11, 1).Value = 0.0008376454934477806
$ws.Cells.Item(11, 2).Value = 0.999739408493042
$ws.Cells.Item(11, 3).Value = 0.002683518454432487
$ws.Cells.Item(11, 4).Value = 0.9995530843734741
$ws.Cells.Item(12, 1).Value = 0.0008292519487440586
$ws.Cells.Item(12, 2).Value = 0.9997193217277527
$ws.Cells.Item(12, 3).Value = 0.005067797843366861
$ws.Cells.Item(12, 4).Value = 0.999106228351593
$ws.Cells.Item(13, 1).Value = 0.0005775628960691392
$ws.Cells.Item(13, 2).Value = 0.9998396039009094
$ws.Cells.Item(13, 3).Value = 0.00634622061625123
$ws.Cells.Item(13, 4).Value = 0.9992551803588867
$ws.Cells.Item(14, 1).Value = 0.0004264594463165849
$ws.Cells.Item(14, 2).Value = 0.9998997449874878
$ws.Cells.Item(14, 3).Value = 0.003960348200052977
$ws.Cells.Item(14, 4).Value = 0.9994041323661804
$ws.Cells.Item(15, 1).Value = 0.0004744645848404616
$ws.Cells.Item(15, 2).Value = 0.9998596906661987
$ws.Cells.Item(15, 3).Value = 0.007094657048583031
$ws.Cells.Item(15, 4).Value = 0.999106228351593
$ws.Cells.Item(16, 1).Value = 0.0003352063067723066
$ws.Cells.Item(16, 2).Value = 0.9999198317527771
$ws.Cells.Item(16, 3).Value = 0.002948341425508261
$ws.Cells.Item(16, 4).Value = 0.9995530843734741
$ws.Cells.Item(17, 1).Value = 0.000135337992105633
$ws.Cells.Item(17, 2).Value = 0.9999799728393555
$ws.Cells.Item(17, 3).Value = 0.003680293913930655
$ws.Cells.Item(17, 4).Value = 0.9995530843734741
$ws.Cells.Item(18, 1).Value = 0.0001632750499993563
$ws.Cells.Item(18, 2).Value = 0.9999398589134216
$ws.Cells.Item(18, 3).Value = 0.005347238853573799
$ws.Cells.Item(18, 4).Value = 0.9995530843734741
$ws.Cells.Item(19, 1).Value = 0.0006788230966776609
$ws.Cells.Item(19, 2).Value = 0.9997594356536865
$ws.Cells.Item(19, 3).Value = 0.01118310820311308
$ws.Cells.Item(19, 4).Value = 0.9992551803588867
$ws.Cells.Item(20, 1).Value = 0.001153093064203858
$ws.Cells.Item(20, 2).Value = 0.9998195767402649
$ws.Cells.Item(20, 3).Value = 0.007868650369346142
$ws.Cells.Item(20, 4).Value = 0.9992551803588867
$ws.Cells.Item(21, 1).Value = 0.0005811135051771998
$ws.Cells.Item(21, 2).Value = 0.9998596906661987
$ws.Cells.Item(21, 3).Value = 0.007795257959514856
$ws.Cells.Item(21, 4).Value = 0.9992551803588867
$ws.Cells.Item(22, 1).Value = 0.0001576094800839201
$ws.Cells.Item(22, 2).Value = 0.9999198317527771
$ws.Cells.Item(22, 3).Value = 0.003422838868573308
$ws.Cells.Item(22, 4).Value = 0.9997020959854126
$ws.Cells.Item(23, 1).Value = 0.0000976195588009432
$ws.Cells.Item(23, 2).Value = 0.9999799728393555
$ws.Cells.Item(23, 3).Value = 0.006421252153813839
$ws.Cells.Item(23, 4).Value = 0.9992551803588867
$ws.Cells.Item(24, 1).Value = 0.0007138215005397797
$ws.Cells.Item(24, 2).Value = 0.9998596906661987
$ws.Cells.Item(24, 3).Value = 0.002709664404392242
$ws.Cells.Item(24, 4).Value = 0.9995530843734741
$ws.Cells.Item(25, 1).Value = 0.0001199529433506541
$ws.Cells.Item(25, 2).Value = 0.9999799728393555
$ws.Cells.Item(25, 3).Value = 0.00490808067843318
$ws.Cells.Item(25, 4).Value = 0.9994041323661804
$ws.Cells.Item(26, 1).Value = 0.0001240165147464722
$ws.Cells.Item(26, 2).Value = 0.9999598860740662
$ws.Cells.Item(26, 3).Value = 0.004653803538531065
$ws.Cells.Item(26, 4).Value = 0.9994041323661804
$ws.Cells.Item(27, 1).Value = 0.0002203768526669592
$ws.Cells.Item(27, 2).Value = 0.9999398589134216
$ws.Cells.Item(27, 3).Value = 0.004867083858698606
$ws.Cells.Item(27, 4).Value = 0.9997020959854126
$ws.Cells.Item(28, 1).Value = 0.0005846436833962798
$ws.Cells.Item(28, 2).Value = 0.9998797178268433
$ws.Cells.Item(28, 3).Value = 0.007863040082156658
$ws.Cells.Item(28, 4).Value = 0.999106228351593
$ws.Cells.Item(29, 1).Value = 0.0006568056414835155
$ws.Cells.Item(29, 2).Value = 0.9998396039009094
$ws.Cells.Item(29, 3).Value = 0.004291488323360682
$ws.Cells.Item(29, 4).Value = 0.9995530843734741
$ws.Cells.Item(30, 1).Value = 0.0001989922748180106
$ws.Cells.Item(30, 2).Value = 0.9999198317527771
$ws.Cells.Item(30, 3).Value = 0.004492656793445349
$ws.Cells.Item(30, 4).Value = 0.9992551803588867
$ws.Cells.Item(31, 1).Value = 0.0001075313048204407
$ws.Cells.Item(31, 2).Value = 0.9999799728393555
$ws.Cells.Item(31, 3).Value = 0.005401885602623224
$ws.Cells.Item(31, 4).Value = 0.9995530843734741
$ws.Cells.Item(32, 1).Value = 0.0002346182736800984
$ws.Cells.Item(32, 2).Value = 0.9999198317527771
$ws.Cells.Item(32, 3).Value = 0.006855674088001251
$ws.Cells.Item(32, 4).Value = 0.9997020959854126
$ws.Cells.Item(33, 1).Value = 0.0001831540430430323
$ws.Cells.Item(33, 2).Value = 0.9999398589134216
$ws.Cells.Item(33, 3).Value = 0.005513347219675779
$ws.Cells.Item(33, 4).Value = 0.9994041323661804
$ws.Cells.Item(34, 1).Value = 0.00007482407090719789
$ws.Cells.Item(34, 2).Value = 0.9999799728393555
$ws.Cells.Item(34, 3).Value = 0.006002883426845074
$ws.Cells.Item(34, 4).Value = 0.9997020959854126
$ws.Cells.Item(35, 1).Value = 0.0002404460683465004
$ws.Cells.Item(35, 2).Value = 0.9999398589134216
$ws.Cells.Item(35, 3).Value = 0.00524580106139183
$ws.Cells.Item(35, 4).Value = 0.9995530843734741
$ws.Cells.Item(36, 1).Value = 0.0002606218622531742
$ws.Cells.Item(36, 2).Value = 0.9998997449874878
$ws.Cells.Item(36, 3).Value = 0.008635876700282097
$ws.Cells.Item(36, 4).Value = 0.9992551803588867
$ws.Cells.Item(37, 1).Value = 0.0003517495933920145
$ws.Cells.Item(37, 2).Value = 0.9999198317527771
$ws.Cells.Item(37, 3).Value = 0.01050938200205564
$ws.Cells.Item(37, 4).Value = 0.9992551803588867
$ws.Cells.Item(38, 1).Value = 0.0006204597884789109
$ws.Cells.Item(38, 2).Value = 0.9998596906661987
$ws.Cells.Item(38, 3).Value = 0.01137394551187754
$ws.Cells.Item(38, 4).Value = 0.999106228351593
$ws.Cells.Item(39, 1).Value = 0.0001711104705464095
$ws.Cells.Item(39, 2).Value = 0.9999799728393555
$ws.Cells.Item(39, 3).Value = 0.004830614197999239
$ws.Cells.Item(39, 4).Value = 0.9995530843734741
$ws.Cells.Item(40, 1).Value = 0.0003749522147700191
$ws.Cells.Item(40, 2).Value = 0.9999198317527771
$ws.Cells.Item(40, 3).Value = 0.002559789922088385
$ws.Cells.Item(40, 4).Value = 0.9997020959854126
$ws.Cells.Item(41, 1).Value = 0.00008093790529528633
$ws.Cells.Item(41, 2).Value = 0.9999799728393555
$ws.Cells.Item(41, 3).Value = 0.004765441175550222
$ws.Cells.Item(41, 4).Value = 0.9997020959854126
$ws.Cells.Item(42, 1).Value = 0.00001200423321279231
$ws.Cells.Item(42, 2).Value = 1
$ws.Cells.Item(42, 3).Value = 0.005059763323515654
$ws.Cells.Item(42, 4).Value = 0.9997020959854126
$ws.Cells.Item(43, 1).Value = 0.00056178227532655
$ws.Cells.Item(43, 2).Value = 0.9998997449874878
$ws.Cells.Item(43, 3).Value = 0.006432718597352505
$ws.Cells.Item(43, 4).Value = 0.9994041323661804
$ws.Cells.Item(44, 1).Value = 0.00001521555077488301
$ws.Cells.Item(44, 2).Value = 1
$ws.Cells.Item(44, 3).Value = 0.007508097682148218
$ws.Cells.Item(44, 4).Value = 0.9994041323661804
$ws.Cells.Item(45, 1).Value = 0.0001871140120783821
$ws.Cells.Item(45, 2).Value = 0.9999398589134216
$ws.Cells.Item(45, 3).Value = 0.006027807481586933
$ws.Cells.Item(45, 4).Value = 0.9997020959854126
$ws.Cells.Item(46, 1).Value = 0.00002265248622279614
$ws.Cells.Item(46, 2).Value = 1
$ws.Cells.Item(46, 3).Value = 0.005880823358893394
$ws.Cells.Item(46, 4).Value = 0.9997020959854126
$ws.Cells.Item(47, 1).Value = 0.0007947962149046361
$ws.Cells.Item(47, 2).Value = 0.9998797178268433
$ws.Cells.Item(47, 3).Value = 0.006195004098117352
$ws.Cells.Item(47, 4).Value = 0.9997020959854126
$ws.Cells.Item(48, 1).Value = 0.0001226189051521942
$ws.Cells.Item(48, 2).Value = 0.9999398589134216
$ws.Cells.Item(48, 3).Value = 0.006707510910928249
$ws.Cells.Item(48, 4).Value = 0.9997020959854126
$ws.Cells.Item(49, 1).Value = 0.0003272075846325606
$ws.Cells.Item(49, 2).Value = 0.9999198317527771
$ws.Cells.Item(49, 3).Value = 0.008327260613441467
$ws.Cells.Item(49, 4).Value = 0.9994041323661804
$ws.Cells.Item(50, 1).Value = 0.0008234516717493534
$ws.Cells.Item(50, 2).Value = 0.9999598860740662
$ws.Cells.Item(50, 3).Value = 0.005264736711978912
$ws.Cells.Item(50, 4).Value = 0.9997020959854126
$ws.Cells.Item(51, 1).Value = 0.000008670182978676166
$ws.Cells.Item(51, 3).Value = 0.005727715324610472
$ws.Cells.Item(51, 4).Value = 0.9997020959854126
